$d = $word.ActiveDocument

# Remove the trailing "Ver no Jupiter..." / "(c) 2020 ..." footer block
# (and the blank paragraph immediately preceding it), which sat right
# after the "Rio de Janeiro: Elsevier Editora, 2007." bibliography line.
$marker = "Rio de Janeiro: Elsevier Editora, 2007."
$footerEnd = "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$paragraphs = $d.Paragraphs
$startPara = $null
$endPara = $null

for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $text = $paragraphs.Item($i).Range.Text
    if ($text -match [regex]::Escape($marker)) {
        # the blank paragraph right after the bibliography entry starts the block to delete
        $startPara = $i + 1
    }
    if ($startPara -ne $null -and $i -ge $startPara -and $text -match [regex]::Escape($footerEnd)) {
        $endPara = $i
        break
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $rangeStart = $paragraphs.Item($startPara).Range.Start
    $rangeEnd = $paragraphs.Item($endPara).Range.End
    $deleteRange = $d.Range($rangeStart, $rangeEnd)
    $deleteRange.Delete()
}
